$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2533654.8
$ws.Range("J70").Value = 4054387.8
$ws.Range("L70").Value = 12163163.4
$ws.Range("N70").Value = -12163703.4
$ws.Range("H73").Value = 2533654.8
$ws.Range("J73").Value = 4054387.8
$ws.Range("L73").Value = 12163163.4
$ws.Range("N73").Value = -12165035.4
$ws.Range("H76").Value = 5238.4
$ws.Range("J76").Value = 5669.8
$ws.Range("L76").Value = 5669.8
$ws.Range("N76").Value = -6299.8
$ws.Range("H79").Value = 5238.4
$ws.Range("J79").Value = 5669.8
$ws.Range("L79").Value = 5669.8
$ws.Range("N79").Value = -7853.8
$ws.Range("H86").Value = 3049.2222
$ws.Range("I86").Value = 2648.4285
$ws.Range("K86").Value = 2648.4285
$ws.Range("M86").Value = -1525.4285
$ws.Range("H89").Value = 3049.2222
$ws.Range("I89").Value = 2648.4285
$ws.Range("K89").Value = 13242.1425
$ws.Range("M89").Value = -7626.1425
$ws.Range("H92").Value = 595.7826
$ws.Range("I92").Value = 473.70587
$ws.Range("K92").Value = 473.70587
$ws.Range("M92").Value = 774.29413
$ws.Range("H94").Value = 3940.1
$ws.Range("I94").Value = 3933.4443
$ws.Range("K94").Value = 3933.4443
$ws.Range("M94").Value = -3482.4443
$ws.Range("H101").Value = 724.25
$ws.Range("I101").Value = 764.6667
$ws.Range("J101").Value = 700
$ws.Range("K101").Value = 2294.0001
$ws.Range("L101").Value = 2100
$ws.Range("M101").Value = -672.0001000000002
$ws.Range("N101").Value = -5344
$ws.Range("H106").Value = 3737.4
$ws.Range("I106").Value = 2496.75
$ws.Range("K106").Value = 2496.75
$ws.Range("M106").Value = -1865.75
$ws.Range("H132").Value = 4224199
$ws.Range("I132").Value = 4352121.5
$ws.Range("K132").Value = 13056364.5
$ws.Range("M132").Value = -13053834.5
$ws.Range("H137").Value = 8804.754999999999
$ws.Range("J137").Value = 3381.7083
$ws.Range("L137").Value = 10145.1249
$ws.Range("N137").Value = -15245.1249

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15129.829
$ws.Range("I32").Value = 15129.829
$ws.Range("K32").Value = 15129.829
$ws.Range("M32").Value = -14842.829
$ws.Range("H61").Value = 3047.8667
$ws.Range("I61").Value = 2123.6365
$ws.Range("K61").Value = 2123.6365
$ws.Range("M61").Value = -1911.6365
$ws.Range("H105").Value = 38185
$ws.Range("J105").Value = 38185
$ws.Range("L105").Value = 38185
$ws.Range("N105").Value = -45173
$ws.Range("H110").Value = 1307.7333
$ws.Range("I110").Value = 1324.2727
$ws.Range("K110").Value = 1324.2727
$ws.Range("M110").Value = 720.7273
$ws.Range("H122").Value = 2014.6154
$ws.Range("I122").Value = 1814.5
$ws.Range("J122").Value = 2681.6667
$ws.Range("K122").Value = 5443.5
$ws.Range("L122").Value = 8045.000100000001
$ws.Range("M122").Value = -2993.5
$ws.Range("N122").Value = -12945.0001
$ws.Range("H132").Value = 1012.86365
$ws.Range("I132").Value = 920.1395
$ws.Range("K132").Value = 2760.4185
$ws.Range("M132").Value = -230.4184999999998
$ws.Range("H133").Value = 67113.664
$ws.Range("J133").Value = 65156.6
$ws.Range("L133").Value = 65156.6
$ws.Range("N133").Value = -70216.60000000001
$ws.Range("H136").Value = 3047.8667
$ws.Range("I136").Value = 2123.6365
$ws.Range("K136").Value = 6370.9095
$ws.Range("M136").Value = -3820.9095

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1149.3704
$ws.Range("I94").Value = 1239.0526
$ws.Range("K94").Value = 1239.0526
$ws.Range("M94").Value = -788.0526
$ws.Range("H105").Value = 2303.8572
$ws.Range("I105").Value = 904.82355
$ws.Range("K105").Value = 904.82355
$ws.Range("M105").Value = 842.17645

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 12722143
$ws.Range("I19").Value = 15900241
$ws.Range("K19").Value = 15900241
$ws.Range("M19").Value = -15900071
$ws.Range("H24").Value = 12722143
$ws.Range("I24").Value = 15900241
$ws.Range("K24").Value = 15900241
$ws.Range("M24").Value = -15900071
$ws.Range("H58").Value = 1721.9688
$ws.Range("I58").Value = 1577.8
$ws.Range("J58").Value = 3884.5
$ws.Range("K58").Value = 1577.8
$ws.Range("L58").Value = 3884.5
$ws.Range("M58").Value = -1374.8
$ws.Range("N58").Value = -4290.5
$ws.Range("H99").Value = 6222.75
$ws.Range("I99").Value = 5210.5
$ws.Range("K99").Value = 5210.5
$ws.Range("M99").Value = -3712.5
$ws.Range("H105").Value = 1809.8667
$ws.Range("I105").Value = 2608.1667
$ws.Range("K105").Value = 2608.1667
$ws.Range("M105").Value = -861.1667000000002
$ws.Range("H122").Value = 11505.19
$ws.Range("I122").Value = 13471.706
$ws.Range("J122").Value = 3147.5
$ws.Range("K122").Value = 40415.118
$ws.Range("L122").Value = 9442.5
$ws.Range("M122").Value = -37965.118
$ws.Range("N122").Value = -14342.5
$ws.Range("H126").Value = 6222.75
$ws.Range("I126").Value = 5210.5
$ws.Range("K126").Value = 15631.5
$ws.Range("M126").Value = -13161.5
$ws.Range("H132").Value = 19109.121
$ws.Range("I132").Value = 19691.031
$ws.Range("K132").Value = 59073.09299999999
$ws.Range("M132").Value = -56543.09299999999
$ws.Range("H134").Value = 1803.909
$ws.Range("I134").Value = 1581.6316
$ws.Range("J134").Value = 3211.6667
$ws.Range("K134").Value = 4744.8948
$ws.Range("L134").Value = 9635.000100000001
$ws.Range("M134").Value = -2209.8948
$ws.Range("N134").Value = -14705.0001
$ws.Range("H136").Value = 1721.9688
$ws.Range("I136").Value = 1577.8
$ws.Range("J136").Value = 3884.5
$ws.Range("K136").Value = 4733.4
$ws.Range("L136").Value = 11653.5
$ws.Range("M136").Value = -2183.4
$ws.Range("N136").Value = -16753.5
$ws.Range("H141").Value = 287999.84
$ws.Range("J141").Value = 287999.84
$ws.Range("L141").Value = 287999.84
$ws.Range("N141").Value = -298359.84

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 593.3333
$ws.Range("I40").Value = 175.33333
$ws.Range("J40").Value = 697.8333
$ws.Range("K40").Value = 701.33332
$ws.Range("L40").Value = 2791.3332
$ws.Range("M40").Value = -632.33332
$ws.Range("N40").Value = -2929.3332
$ws.Range("H92").Value = 1244.4667
$ws.Range("I92").Value = 2175.6428
$ws.Range("J92").Value = 429.6875
$ws.Range("K92").Value = 6526.928400000001
$ws.Range("L92").Value = 1289.0625
$ws.Range("M92").Value = -5278.928400000001
$ws.Range("N92").Value = -3785.0625
$ws.Range("H122").Value = 956.4516
$ws.Range("J122").Value = 1083.409
$ws.Range("L122").Value = 9750.681
$ws.Range("N122").Value = -14650.681
$ws.Range("H132").Value = 2056.75
$ws.Range("J132").Value = 1993.6
$ws.Range("L132").Value = 17942.4
$ws.Range("N132").Value = -23002.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5852.68
$ws.Range("I80").Value = 3911.4211
$ws.Range("K80").Value = 3911.4211
$ws.Range("M80").Value = -2913.4211
$ws.Range("H83").Value = 5852.68
$ws.Range("I83").Value = 3911.4211
$ws.Range("K83").Value = 19557.1055
$ws.Range("M83").Value = -14565.1055
$ws.Range("H97").Value = 952.4727
$ws.Range("I97").Value = 836.9459000000001
$ws.Range("J97").Value = 1189.9445
$ws.Range("K97").Value = 836.9459000000001
$ws.Range("L97").Value = 1189.9445
$ws.Range("M97").Value = -340.9459000000001
$ws.Range("N97").Value = -2181.9445

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 892.0952
$ws.Range("I22").Value = 633.1667
$ws.Range("J22").Value = 1237.3334
$ws.Range("K22").Value = 633.1667
$ws.Range("L22").Value = 1237.3334
$ws.Range("M22").Value = -338.1667
$ws.Range("N22").Value = -1827.3334
$ws.Range("H27").Value = 892.0952
$ws.Range("I27").Value = 633.1667
$ws.Range("J27").Value = 1237.3334
$ws.Range("K27").Value = 633.1667
$ws.Range("L27").Value = 1237.3334
$ws.Range("M27").Value = -526.1667
$ws.Range("N27").Value = -1451.3334
$ws.Range("H46").Value = 4836.5
$ws.Range("I46").Value = 659.75
$ws.Range("J46").Value = 8177.9
$ws.Range("K46").Value = 659.75
$ws.Range("L46").Value = 8177.9
$ws.Range("M46").Value = -471.75
$ws.Range("N46").Value = -8553.9

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6482.75
$ws.Range("I62").Value = 5793.8
$ws.Range("J62").Value = 6974.857
$ws.Range("K62").Value = 5793.8
$ws.Range("L62").Value = 6974.857
$ws.Range("M62").Value = -5169.8
$ws.Range("N62").Value = -8222.857
$ws.Range("H65").Value = 6482.75
$ws.Range("I65").Value = 5793.8
$ws.Range("J65").Value = 6974.857
$ws.Range("K65").Value = 28969
$ws.Range("L65").Value = 34874.285
$ws.Range("M65").Value = -25849
$ws.Range("N65").Value = -41114.285
$ws.Range("H126").Value = 314524.56
$ws.Range("I126").Value = 2131.9092
$ws.Range("K126").Value = 6395.7276
$ws.Range("M126").Value = -3925.7276
